$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: "altitude" -> "height_above_snow_surface" ---------------------
# Rename the variable name in column A.
$ws.Range("A2").Value = "height_above_snow_surface"

# The "standard_name" row (row 6: standard_name / altitude) is removed because
# height_above_snow_surface has no CF standard name.
$ws.Rows(6).Delete()

# valid_min / valid_max (now rows 9 & 10 after the delete) become derived values
$ws.Range("C9").Value = "<derived>"
$ws.Range("C10").Value = "<derived>"

# Insert a new "cell_methods" / "time: mean" row before "coordinates" (row 11)
$ws.Rows(11).Insert()
$ws.Range("B11").Value = "cell_methods"
$ws.Range("C11").Value = "time: mean"
# Fix up formatting for the freshly-inserted row: B11 plain/default, C11 like C15 (style 9)
$ws.Range("B11").ClearFormats()
$ws.Range("C15").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# --- Block 2: wind_speed valid_min / valid_max become derived --------------
$ws.Range("C21").Value = "<derived>"
$ws.Range("C22").Value = "<derived>"

# --- Block 3: wind_from_direction valid_min / valid_max become derived -----
$ws.Range("C33").Value = "<derived>"
$ws.Range("C34").Value = "<derived>"

# --- Block 4: qc_flag - remove blank standard_name row (row 42) ------------
$ws.Rows(42).Delete()

# --- Column A width ----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 38.92

# --- Sheet view: drop frozen/scrolled topLeftCell, update active selection -
$ws.Range("C45").Select()
